$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 34 ("Agrícola del
# Norte S.A. de Arica" - Perejil weekly update). Everything currently at
# row 34 and below shifts down by two rows, matching the target diff where
# old row 34 -> new row 36, ..., old row 59 -> new row 61.
$ws.Rows("34:35").Insert()

# New row 34: Perejil, "Primera" quality, week of 2023-10-26 (serial 45225)
$ws.Cells.Item(34, 1).Value  = 1
$ws.Cells.Item(34, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(34, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(34, 4).Value  = 45225
$ws.Cells.Item(34, 5).Value  = 15
$ws.Cells.Item(34, 6).Value  = 100112044
$ws.Cells.Item(34, 7).Value  = "Perejil"
$ws.Cells.Item(34, 8).Value  = "Sin especificar"
$ws.Cells.Item(34, 9).Value  = "Primera"
$ws.Cells.Item(34, 10).Value = 250
$ws.Cells.Item(34, 11).Value = 2000
$ws.Cells.Item(34, 12).Value = 2000
$ws.Cells.Item(34, 13).Value = 2000
$ws.Cells.Item(34, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(34, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(34, 16).Value = 1000
$ws.Cells.Item(34, 17).Value = 2
$ws.Cells.Item(34, 18).Value = "Hortaliza"

# New row 35: Perejil, "Segunda" quality, same week (serial 45225)
$ws.Cells.Item(35, 1).Value  = 1
$ws.Cells.Item(35, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(35, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(35, 4).Value  = 45225
$ws.Cells.Item(35, 5).Value  = 15
$ws.Cells.Item(35, 6).Value  = 100112044
$ws.Cells.Item(35, 7).Value  = "Perejil"
$ws.Cells.Item(35, 8).Value  = "Sin especificar"
$ws.Cells.Item(35, 9).Value  = "Segunda"
$ws.Cells.Item(35, 10).Value = 110
$ws.Cells.Item(35, 11).Value = 1700
$ws.Cells.Item(35, 12).Value = 1700
$ws.Cells.Item(35, 13).Value = 1700
$ws.Cells.Item(35, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(35, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(35, 16).Value = 850
$ws.Cells.Item(35, 17).Value = 2
$ws.Cells.Item(35, 18).Value = "Hortaliza"
